# Add a new data row (row 63) to Sheet1, continuing the daily log table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 63

# Column A holds a date-like string ("yyyy/mm/dd") that must stay plain
# text, matching every other row in the table. Mark it as text first so
# Excel doesn't auto-convert it into a real date value/format, then
# restore the default "Normal" style so no extra per-cell formatting is
# left behind (the other data rows carry no explicit style either).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/05"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "日"
$ws.Cells.Item($row, 3).Value = 4
$ws.Cells.Item($row, 4).Value = 201
